# Apply crypto price/volume update as described in commit
# "Updated cryptos list on Mon Nov 27 19:30:05 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.021.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.007.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.305.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.740"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.002.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.940.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0816"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  -5.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.15%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.124"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0609"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.473.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.79%  "
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.193.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.96%  "
